$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7080510854721069
$ws.Range("B1").Value = 1.074634671211243
$ws.Range("C1").Value = 2.353059053421021
$ws.Range("D1").Value = 3.628820896148682
$ws.Range("E1").Value = 1.697909593582153
